$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2
Set-TextValue "D2" "61.439.39"
Set-TextValue "E2" "  +0.86%  "

# Row 3
Set-TextValue "D3" "3.393.33"
Set-TextValue "E3" "  +0.13%  "

# Row 5
Set-TextValue "D5" "577.22"
Set-TextValue "E5" "  +1.10%  "

# Row 6
Set-TextValue "D6" "141.09"
Set-TextValue "E6" "  -0.50%  "

# Row 8
Set-TextValue "E8" "  -0.28%  "

# Row 9
Set-TextValue "D9" "7.70"
Set-TextValue "E9" "  +2.08%  "

# Row 10
Set-TextValue "E10" "  -0.57%  "

# Row 11
Set-TextValue "E11" "  -1.57%  "

# Row 12
Set-TextValue "D12" "3.972.78"
Set-TextValue "E12" "  +0.09%  "

# Row 13
Set-TextValue "E13" "  +0.46%  "

# Row 14
Set-TextValue "D14" "28.28"
Set-TextValue "E14" "  -0.11%  "

# Row 15
Set-TextValue "B15" "ShibaInu"
Set-TextValue "C15" "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D15" "0.0000171"
Set-TextValue "E15" "  +0.29%  "

# Row 16
Set-TextValue "B16" "WrappedEther"
Set-TextValue "C16" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D16" "3.388.90"
Set-TextValue "E16" "  -0.29%  "

# Row 17
Set-TextValue "D17" "61.431.00"
Set-TextValue "E17" "  +0.75%  "

# Row 18
Set-TextValue "D18" "6.14"
Set-TextValue "E18" "  -0.56%  "

# Row 19
Set-TextValue "D19" "13.67"
Set-TextValue "E19" "  -1.51%  "

# Row 20
Set-TextValue "D20" "8.99"
Set-TextValue "E20" "  +0.03%  "

# Row 21
Set-TextValue "D21" "391.75"
Set-TextValue "E21" "  +1.78%  "

# Row 22
Set-TextValue "D22" "75.36"
Set-TextValue "E22" "  +1.59%  "

# Row 23
Set-TextValue "E23" "  -0.42%  "

# Row 24
Set-TextValue "E24" "  +0.01%  "

# Row 25
Set-TextValue "D25" "0.0000113"
Set-TextValue "E25" "  -3.32%  "

# Row 26
Set-TextValue "D26" "0.192"
Set-TextValue "E26" "  +7.81%  "

# Row 27
Set-TextValue "E27" "  -0.05%  "

# Row 28
Set-TextValue "D28" "7.28"
Set-TextValue "E28" "  -1.45%  "

# Row 29
Set-TextValue "E29" "  +1.10%  "

# Row 30
Set-TextValue "E30" "  +0.62%  "

# Row 31
Set-TextValue "E31" "  -0.03%  "

# Row 32
Set-TextValue "E32" "  -4.43%  "

# Row 33
Set-TextValue "D33" "23.41"
Set-TextValue "E33" "  -0.36%  "

# Row 34
Set-TextValue "E34" "  -0.53%  "

# Row 35
Set-TextValue "D35" "167.78"
Set-TextValue "E35" "  +0.26%  "

# Row 36
Set-TextValue "D36" "5.05"
Set-TextValue "E36" "  +1.42%  "

# Row 37
Set-TextValue "D37" "3.424.93"
Set-TextValue "E37" "  +0.13%  "

# Row 38
Set-TextValue "E38" "  -0.39%  "

# Row 39
Set-TextValue "D39" "0.0769"
Set-TextValue "E39" "  -0.63%  "

# Row 40
Set-TextValue "D40" "26.27"
Set-TextValue "E40" "  -4.26%  "

# Row 41
Set-TextValue "E41" "  -0.07%  "

# Row 42
Set-TextValue "D42" "4.43"
Set-TextValue "E42" "  +0.20%  "

# Row 43
Set-TextValue "E43" "  -0.62%  "

# Row 44
Set-TextValue "E44" "  +1.05%  "

# Row 45
Set-TextValue "D45" "2.465.86"
Set-TextValue "E45" "  -0.82%  "

# Row 46
Set-TextValue "D46" "22.92"
Set-TextValue "E46" "  -0.26%  "

# Row 47
Set-TextValue "E47" "  -1.63%  "

# Row 48
Set-TextValue "D48" "1.00"
Set-TextValue "E48" "  +0.09%  "

# Row 49
Set-TextValue "E49" "  -2.11%  "

# Row 50
Set-TextValue "D50" "2.07"
Set-TextValue "E50" "  -0.90%  "

# Row 51
Set-TextValue "E51" "  -1.31%  "
